# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q3").
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $zj)
$q4.Name = "2022-Q4"

# Fetch sheet references AFTER the insert/rename above so indices are fresh.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Copy header-row formatting (bold + border style) from the 2022-Q3 sheet so
# the new sheet matches the existing look (style index carried via format copy).
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "519019"
$q4.Range("B2").Style = "Normal"
$q4.Range("C2").Value = "大成景阳领先混合"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "20.42"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "90.28"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "4.58"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.9352"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 5

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "090016"
$q4.Range("B3").Style = "Normal"
$q4.Range("C3").Value = "大成消费主题混合"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "5.72"
$q4.Range("D3").Style = "Normal"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "91.78"
$q4.Range("E3").Style = "Normal"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "4.83"
$q4.Range("F3").Style = "Normal"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.2763"
$q4.Range("G3").Style = "Normal"
$q4.Range("H3").Value = 5

# Row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "003956"
$q4.Range("B4").Style = "Normal"
$q4.Range("C4").Value = "南方产业智选股票"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "3.60"
$q4.Range("D4").Style = "Normal"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "85.80"
$q4.Range("E4").Style = "Normal"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "4.67"
$q4.Range("F4").Style = "Normal"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.1681"
$q4.Range("G4").Style = "Normal"
$q4.Range("H4").Value = 7

# Row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").NumberFormat = "@"
$q4.Range("B5").Value = "002319"
$q4.Range("B5").Style = "Normal"
$q4.Range("C5").Value = "大成一带一路灵活配置混合"
$q4.Range("D5").NumberFormat = "@"
$q4.Range("D5").Value = "1.49"
$q4.Range("D5").Style = "Normal"
$q4.Range("E5").NumberFormat = "@"
$q4.Range("E5").Value = "92.29"
$q4.Range("E5").Style = "Normal"
$q4.Range("F5").NumberFormat = "@"
$q4.Range("F5").Value = "5.56"
$q4.Range("F5").Style = "Normal"
$q4.Range("G5").NumberFormat = "@"
$q4.Range("G5").Value = "0.0828"
$q4.Range("G5").Style = "Normal"
$q4.Range("H5").Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push existing data rows down by one
#    and insert the new 2022-Q4 summary row at the top (row 2).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

# shift B/C/D values down a row (A column index stays 0,1,2 and gets a new 3)
$ws.Range("B5").Value = $ws.Range("B4").Value()
$ws.Range("C5").Value = $ws.Range("C4").Value()
$ws.Range("D5").Value = $ws.Range("D4").Value()

$ws.Range("B4").Value = $ws.Range("B3").Value()
$ws.Range("C4").Value = $ws.Range("C3").Value()
$ws.Range("D4").Value = $ws.Range("D3").Value()

$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()

$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 1.46

# new row 5 needs the A-column index style/value (copy format from A2, which
# already carries the bold+border "index column" style)
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 3) Restore the originally-active sheet/tab (adding a sheet makes it active
#    by default, but "2021-Q1" was the selected tab before this edit).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
